# Replace the intro paragraph of the Cygnus Finnish activity guide with the
# new wording that names the Cygnus constellation, collapsing the many
# runs that make up the old paragraph into a single plain run.

$d = $word.ActiveDocument

$newText = "Osallistut maailmanlaajuiseen kampanjaan tarkkaillaksesi ja tallentaaksesi himmeimpiä näkyvissä olevia tähtiä keinona mitata valonsaastetta tietyssä paikassa. Paikallistamalla ja tarkkailemalla Cygnus-tähdistö miten valosaaste syntyy kunkin taajaman tai muun ihmisen toiminnan valoista. Antamasi tiedot päivittyvät heti verkossa olevaan tietokantaan, ja näin saadaan käsitys siitä minkä verran taivaan tähdistä on missäkin nähtävissä."

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "Osallistut maailmanlaajuiseen*") {
        # Range covering the paragraph's content, excluding the trailing
        # paragraph mark, so the paragraph itself is preserved.
        $content = $d.Range($r.Start, $r.End - 1)
        $content.Delete()

        # Insert fresh text into the now-empty paragraph; this produces a
        # single run with no run properties (matches a freshly typed run).
        $insertionPoint = $d.Range($r.Start, $r.Start)
        $insertionPoint.InsertBefore($newText)
        break
    }
}
